$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = 0.08059215988451403
$ws.Range("J2").Value = 0.08059215988451404
$ws.Range("M2").Value = 0.74396
$ws.Range("N2").Value = 2.23188
$ws.Range("O2").Value = 0.006259003216804254
$ws.Range("P2").Value = 0.006259003216804255
$ws.Range("Q2").Value = 0.14262754744
$ws.Range("R2").Value = 1.28364792696
$ws.Range("S2").Value = 0.0005044265879663761
$ws.Range("T2").Value = 0.0005044265879663762

# Row 3
$ws.Range("I3").Value = 0.08059215988451403
$ws.Range("J3").Value = 0.08059215988451404
$ws.Range("M3").Value = 88.14978533333333
$ws.Range("O3").Value = 0.7416121699579786
$ws.Range("P3").Value = 0.7416121699579786
$ws.Range("S3").Value = 0.05976812657355481
$ws.Range("T3").Value = 0.05976812657355481

# Row 4
$ws.Range("I4").Value = 0.08059215988451403
$ws.Range("J4").Value = 0.08059215988451404
$ws.Range("M4").Value = 29.76859933333333
$ws.Range("N4").Value = 89.305798
$ws.Range("O4").Value = 0.2504459365921425
$ws.Range("P4").Value = 0.2504459365921425
$ws.Range("Q4").Value = 5.707057252590666
$ws.Range("R4").Value = 51.363515273316
$ws.Range("S4").Value = 0.02018397896426081
$ws.Range("T4").Value = 0.02018397896426081

# Row 5
$ws.Range("I5").Value = 0.08059215988451403
$ws.Range("J5").Value = 0.08059215988451404
$ws.Range("M5").Value = 0.2000323333333334
$ws.Range("N5").Value = 0.6000970000000001
$ws.Range("O5").Value = 0.00168289023307462
$ws.Range("P5").Value = 0.00168289023307462
$ws.Range("Q5").Value = 0.03834899875266668
$ws.Range("R5").Value = 0.3451409887740001
$ws.Range("S5").Value = 0.0001356277587320369
$ws.Range("T5").Value = 0.0001356277587320369

# Row 6
$ws.Range("G6").Value = 2.187103
$ws.Range("H6").Value = 6.561309
$ws.Range("I6").Value = 0.9194078401154859
$ws.Range("J6").Value = 0.919407840115486
$ws.Range("M6").Value = 0.74396
$ws.Range("N6").Value = 2.23188
$ws.Range("O6").Value = 0.006259003216804254
$ws.Range("P6").Value = 0.006259003216804255
$ws.Range("Q6").Value = 1.62711714788
$ws.Range("R6").Value = 14.64405433092
$ws.Range("S6").Value = 0.005754576628837878
$ws.Range("T6").Value = 0.005754576628837879

# Row 7
$ws.Range("G7").Value = 2.187103
$ws.Range("H7").Value = 6.561309
$ws.Range("I7").Value = 0.9194078401154859
$ws.Range("J7").Value = 0.919407840115486
$ws.Range("M7").Value = 88.14978533333333
$ws.Range("O7").Value = 0.7416121699579786
$ws.Range("P7").Value = 0.7416121699579786
$ws.Range("Q7").Value = 192.7926599518893
$ws.Range("R7").Value = 1735.133939567004
$ws.Range("S7").Value = 0.6818440433844238
$ws.Range("T7").Value = 0.6818440433844238

# Row 8
$ws.Range("G8").Value = 2.187103
$ws.Range("H8").Value = 6.561309
$ws.Range("I8").Value = 0.9194078401154859
$ws.Range("J8").Value = 0.919407840115486
$ws.Range("M8").Value = 29.76859933333333
$ws.Range("N8").Value = 89.305798
$ws.Range("O8").Value = 0.2504459365921425
$ws.Range("P8").Value = 0.2504459365921425
$ws.Range("Q8").Value = 65.10699290773132
$ws.Range("R8").Value = 585.9629361695819
$ws.Range("S8").Value = 0.2302619576278817
$ws.Range("T8").Value = 0.2302619576278817

# Row 9
$ws.Range("G9").Value = 2.187103
$ws.Range("H9").Value = 6.561309
$ws.Range("I9").Value = 0.9194078401154859
$ws.Range("J9").Value = 0.919407840115486
$ws.Range("M9").Value = 0.2000323333333334
$ws.Range("N9").Value = 0.6000970000000001
$ws.Range("O9").Value = 0.00168289023307462
$ws.Range("P9").Value = 0.00168289023307462
$ws.Range("Q9").Value = 0.4374913163303334
$ws.Range("R9").Value = 3.937421846973
$ws.Range("S9").Value = 0.001547262474342583
$ws.Range("T9").Value = 0.001547262474342584
